$wb = $excel.ActiveWorkbook

# --- Sheet "Observations": append the new Add Debit/Credit Cards test row ---
$obs = $wb.Worksheets.Item("Observations")

$obs.Range("A4").Value = 3

# Same date format as the rows above it (d-mmm-yy) so it reuses that style.
$obs.Range("B4").Value = 45142
$obs.Range("B4").NumberFormat = "d-mmm-yy"

$obs.Range("C4").Value = "Add Debit/Credit Cards"
$obs.Range("D4").Value = "Add Debit/Credit Cards"
$obs.Range("F4").Value = "Trying to add the existing card,its refreshing the card fields, after entering the valid new card details, address is not auto populating"
$obs.Range("E4").Value = "Billing Address"

# Widen column C so the new "Add Debit/Credit Cards" text fits (was grouped
# with column D at 15.5 before; now split out on its own).
$obs.Columns.Item(3).ColumnWidth = 23.75

# Move the active selection to the newly entered cell.
$obs.Range("E4").Select() | Out-Null

# --- Sheet "Blockers": move the active selection ---
$blk = $wb.Worksheets.Item("Blockers")
$blk.Range("C3").Select() | Out-Null

$obs.Activate() | Out-Null
